$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Dead volume" row of inputs (row 6): label, value, unit
$ws.Range("A6").Value = "Dead volume"
$ws.Range("B6").Value = 0.75
$ws.Range("C6").Value = "mL"

# Match the yellow input-cell look (same fill as the other input cells)
# plus right-aligned text for the new value cell.
$ws.Range("B6").Interior.Color = 65535
$ws.Range("B6").HorizontalAlignment = -4152

# Update the selection / scroll position that Excel persists with the sheet.
$ws.Range("B6").Select()

# Restore the window to a maximized-looking state (matches the saved view).
$win = $excel.ActiveWindow
$win.WindowState = -4137
